# Update quantities in column C: every 75 becomes 100, every 100 becomes 150.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 208; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2
    if ($v -eq 75) {
        $cell.Value = 100
    } elseif ($v -eq 100) {
        $cell.Value = 150
    }
}

# Record a (hidden, sheet-scoped) AutoFilter database name, as Excel does
# when a filter range was defined on the sheet.
$name = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet2!`$A`$1:`$C`$208")
$name.Visible = $false

# Move the active selection to A2.
$ws.Range("A2").Select()
